# Generate Report for Handback
#
# The "3aeeb515-b459-4b03-83f7-91518af0c92c" row (row 7) on both the
# zh-cn and de-de sheets gets a handback result recorded: the handback
# was detected, but the handback file version is stale, so:
#   - "Latest Target File" (I) gets filled in with the source file name,
#     hyperlinked back to the source markdown file on GitHub (same
#     target as the "Source File Name" (A) link already on that row).
#   - "Latest Handback File" (J) gets filled in with the xlf file name
#     that was handed back.
#   - "Latest Handback DateTime" (K) records when the (stale) handback
#     was processed.
#   - "Error Detail" (P) explains that the handback isn't based on the
#     latest source commit.

$wb = $excel.ActiveWorkbook

$sourceUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/054d6cd52d524b1fea1a1632094bf33ce8f24459/e2e/3aeeb515-b459-4b03-83f7-91518af0c92c.md"
$sourceDisplay = "3aeeb515-b459-4b03-83f7-91518af0c92c.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48e0e421cbce1bbcd2c4768ef1234af9fef91241/e2e/3aeeb515-b459-4b03-83f7-91518af0c92c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/054d6cd52d524b1fea1a1632094bf33ce8f24459/e2e/3aeeb515-b459-4b03-83f7-91518af0c92c.md."

# --- zh-cn sheet, row 7 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I7"), $sourceUrl, "", "", $sourceDisplay)
$wsZhCn.Range("J7").Value = "3aeeb515-b459-4b03-83f7-91518af0c92c.7b929fc771d854f42c0405cfe192faa645138d42.zh-cn.xlf"
$wsZhCn.Range("K7").Value = "2016-08-17 20:56:14"
$wsZhCn.Range("P7").Value = $errorDetail

# --- de-de sheet, row 7 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I7"), $sourceUrl, "", "", $sourceDisplay)
$wsDeDe.Range("J7").Value = "3aeeb515-b459-4b03-83f7-91518af0c92c.7b929fc771d854f42c0405cfe192faa645138d42.de-de.xlf"
$wsDeDe.Range("K7").Value = "2016-08-17 20:56:22"
$wsDeDe.Range("P7").Value = $errorDetail
